$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source workbook stores Price/Volume figures as plain text (inline strings),
# even when they look like plain decimals (e.g. "561.73"). Excel normally auto-
# detects such strings as numbers when assigned via .Value, which would both change
# the cell type and introduce floating-point artifacts (e.g. "0.0000230" -> 2.3E-05).
# Setting NumberFormat to "@" (Text) on each target cell before writing its value
# keeps the data as literal text, matching the source.

$ws.Range("D2").Value = "59.332.33"
$ws.Range("E2").Value = "  +2.18%  "
$ws.Range("D3").Value = "2.993.25"
$ws.Range("E3").Value = "  +0.70%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "561.73"
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.40"
$ws.Range("E6").Value = "  +3.47%  "
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("E8").Value = "  +0.94%  "
$ws.Range("D9").Value = "2.981.74"
$ws.Range("E9").Value = "  +0.65%  "
$ws.Range("E10").Value = "  +2.52%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.17"
$ws.Range("E11").Value = "  +5.83%  "
$ws.Range("E12").Value = "  +2.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000230"
$ws.Range("E13").Value = "  +2.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.72"
$ws.Range("E14").Value = "  +1.87%  "
$ws.Range("E15").Value = "  +1.69%  "
$ws.Range("D16").Value = "3.487.80"
$ws.Range("E16").Value = "  +0.76%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.30"
$ws.Range("E17").Value = "  +6.78%  "
$ws.Range("D18").Value = "2.993.68"
$ws.Range("E18").Value = "  +0.80%  "
$ws.Range("D19").Value = "59.358.05"
$ws.Range("E19").Value = "  +2.42%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "430.54"
$ws.Range("E20").Value = "  +2.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.60"
$ws.Range("E21").Value = "  +2.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.720"
$ws.Range("E22").Value = "  +4.52%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.52"
$ws.Range("E23").Value = "  +3.29%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.13"
$ws.Range("E24").Value = "  +1.58%  "
$ws.Range("E25").Value = "  +0.80%  "
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.23"
$ws.Range("E27").Value = "  +9.97%  "
$ws.Range("E28").Value = "  +0.23%  "
$ws.Range("E29").Value = "  +1.43%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.87"
$ws.Range("E30").Value = "  +3.23%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.74"
$ws.Range("E31").Value = "  +1.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.11"
$ws.Range("E32").Value = "  -0.06%  "
$ws.Range("E33").Value = "  +1.65%  "
$ws.Range("E34").Value = "  +5.84%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.96"
$ws.Range("E35").Value = "  +4.90%  "
$ws.Range("D36").Value = "0.0₃0763"
$ws.Range("E36").Value = "  +9.14%  "
$ws.Range("E37").Value = "  -1.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.86"
$ws.Range("E38").Value = "  +0.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.68"
$ws.Range("E39").Value = "  -0.38%  "
$ws.Range("E40").Value = "  +5.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "407.33"
$ws.Range("E41").Value = "  +7.32%  "
$ws.Range("E42").Value = "  +0.35%  "
$ws.Range("D43").Value = "2.774.73"
$ws.Range("E43").Value = "  +2.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.107"
$ws.Range("E44").Value = "  -1.25%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.252"
$ws.Range("E45").Value = "  +4.05%  "
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "123.44"
$ws.Range("E47").Value = "  +0.60%  "
$ws.Range("B48").Value = "Arweave"
$ws.Range("C48").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "34.63"
$ws.Range("E48").Value = "  +20.40%  "
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.01"
$ws.Range("E50").Value = "  +0.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.52"
$ws.Range("E51").Value = "  -0.44%  "
